# Update cryptos price (D) and 1h volume/change (E) columns
# to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column prices are text that often LOOKS numeric (e.g. "229.73").
# A plain Range.Value assignment would make Excel auto-convert such
# strings into real numbers, which does not match the source data
# (the workbook stores every Price cell as text/inline string).
# Forcing a Text number format before the write keeps the assignment
# a string, and restoring the "Normal" style afterwards avoids leaving
# a stray text-format style on the cell.
function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "38.122.92"
$ws.Range("E2").Value = "  +2.84%  "

# Row 3
Set-TextValue "D3" "2.056.78"
$ws.Range("E3").Value = "  +2.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue "D5" "229.73"
$ws.Range("E5").Value = "  +1.57%  "

# Row 6
$ws.Range("E6").Value = "  +1.89%  "

# Row 7
Set-TextValue "D7" "59.20"
$ws.Range("E7").Value = "  +7.75%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
Set-TextValue "D9" "0.387"
$ws.Range("E9").Value = "  +3.23%  "

# Row 10
$ws.Range("E10").Value = "  +3.96%  "

# Row 11
$ws.Range("E11").Value = "  +2.23%  "

# Row 12
Set-TextValue "D12" "2.359.40"
$ws.Range("E12").Value = "  +2.12%  "

# Row 13
Set-TextValue "D13" "14.69"
$ws.Range("E13").Value = "  +4.63%  "

# Row 14
Set-TextValue "D14" "21.02"
$ws.Range("E14").Value = "  +6.16%  "

# Row 15
Set-TextValue "D15" "0.754"
$ws.Range("E15").Value = "  +2.28%  "

# Row 16
Set-TextValue "D16" "5.29"
$ws.Range("E16").Value = "  +1.63%  "

# Row 17
Set-TextValue "D17" "2.079.53"
$ws.Range("E17").Value = "  +3.33%  "

# Row 18
Set-TextValue "D18" "38.011.87"
$ws.Range("E18").Value = "  +2.77%  "

# Row 19
Set-TextValue "D19" "6.31"
$ws.Range("E19").Value = "  +0.77%  "

# Row 20
Set-TextValue "D20" "69.86"
$ws.Range("E20").Value = "  +2.38%  "

# Row 21
$ws.Range("E21").Value = "  +2.79%  "

# Row 22
Set-TextValue "D22" "224.65"
$ws.Range("E22").Value = "  +0.86%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("E25").Value = "  +3.70%  "

# Row 26
Set-TextValue "D26" "9.29"
$ws.Range("E26").Value = "  +3.39%  "

# Row 27
Set-TextValue "D27" "166.38"
$ws.Range("E27").Value = "  +1.11%  "

# Row 28
Set-TextValue "D28" "0.133"
$ws.Range("E28").Value = "  +6.94%  "

# Row 29
Set-TextValue "D29" "19.05"
$ws.Range("E29").Value = "  +2.62%  "

# Row 30
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
Set-TextValue "D31" "0.120"
$ws.Range("E31").Value = "  +2.48%  "

# Row 32
Set-TextValue "D32" "4.55"
$ws.Range("E32").Value = "  +2.02%  "

# Row 33
Set-TextValue "D33" "4.60"
$ws.Range("E33").Value = "  +2.43%  "

# Row 34
Set-TextValue "D34" "2.07"
$ws.Range("E34").Value = "  +10.88%  "

# Row 35
$ws.Range("E35").Value = "  +0.99%  "

# Row 36
$ws.Range("E36").Value = "  +0.09%  "

# Row 37
Set-TextValue "D37" "6.09"
$ws.Range("E37").Value = "  +14.10%  "

# Row 38
$ws.Range("E38").Value = "  +5.12%  "

# Row 39
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
Set-TextValue "D40" "1.535.60"
$ws.Range("E40").Value = "  +5.46%  "

# Row 41
Set-TextValue "D41" "98.28"
$ws.Range("E41").Value = "  +3.60%  "

# Row 42
$ws.Range("E42").Value = "  +2.29%  "

# Row 43
$ws.Range("E43").Value = "  +4.13%  "

# Row 44
Set-TextValue "D44" "16.83"
$ws.Range("E44").Value = "  +5.74%  "

# Row 45
Set-TextValue "D45" "0.0924"
$ws.Range("E45").Value = "  +2.08%  "

# Row 46
Set-TextValue "D46" "1.14"
$ws.Range("E46").Value = "  +1.14%  "

# Row 47
$ws.Range("E47").Value = "  +13.65%  "

# Row 49
$ws.Range("E49").Value = "  +2.56%  "

# Row 50
Set-TextValue "D50" "7.14"
$ws.Range("E50").Value = "  +0.26%  "

# Row 51
Set-TextValue "D51" "2.247.15"
$ws.Range("E51").Value = "  +2.36%  "
